$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.351.76'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.867.45'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.43'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4708'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2869'
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06568'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.46'
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07873'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.96'
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '1.870.39'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6911'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.108'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.47'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '30.313.58'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.95'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("E19").Value = '  +3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '2.113.34'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.186'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.407'
$ws.Range("E25").Value = '  +2.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.33'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.88'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.362'
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09921'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.377'
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.460'
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04748'
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7031'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.790'
$ws.Range("E39").Value = '  +6.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.295'
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.41'
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.949'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8428'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.99'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '977.12'
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.109'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.139'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.51'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05677'
$ws.Range("E51").Value = '  +0.48%  '
